$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D-column price cells that would otherwise be
# auto-converted to floating point numbers by Excel (values keep trailing
# zeros / exact decimal text in the source workbook).
$textCells = @("D5","D6","D8","D10","D13","D18","D19","D20","D21","D23","D28","D35","D39","D40","D41","D44","D45","D46","D47","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cryptos list values
$ws.Range("D2").Value = "58.074.38"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "2.360.07"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "541.72"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").Value = "136.09"
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "0.561"
$ws.Range("E8").Value = "  +5.15%  "
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").Value = "5.59"
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "23.99"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "2.780.16"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "58.047.69"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "2.359.54"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "10.76"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").Value = "333.03"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "6.81"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "62.85"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("D28").Value = "172.88"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").Value = "0.0₃0741"
$ws.Range("E30").Value = "  +2.56%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  +11.92%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "4.25"
$ws.Range("E35").Value = "  +6.89%  "
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("D39").Value = "39.41"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "145.52"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("D41").Value = "294.03"
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("D44").Value = "0.0947"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("D45").Value = "19.28"
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").Value = "0.0504"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").Value = "0.565"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("D49").Value = "17.54"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "11.06"
$ws.Range("E51").Value = "  +0.40%  "
